$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same (reused) style index instead of
# minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the new data cells
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
